# "Generate Report for Archive"
#
# 1) The localization status for both in-flight docs moves from
#    "Ready for handoff" to "In Translation" - update every cell that
#    carries that status (the Overview rollup sheet's zh-cn/de-de columns,
#    plus each language sheet's own Status column).
# 2) The Status-ish columns that used to be sized for the longer
#    "Ready for handoff" label are narrowed to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Narrow the now-shorter status columns to fit "In Translation".
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
